$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.921.65"
$ws.Range("E2").Value = "  -0.87%  "
$ws.Range("D3").Value = "1.817.54"
$ws.Range("E3").Value = "  +0.28%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.25"
$ws.Range("E5").Value = "  -0.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4658"
$ws.Range("E7").Value = "  +1.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3700"
$ws.Range("E8").Value = "  -1.20%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07379"
$ws.Range("E9").Value = "  -0.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8718"
$ws.Range("E10").Value = "  +0.92%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.45"
$ws.Range("E11").Value = "  -0.60%  "
$ws.Range("D12").Value = "1.847.28"
$ws.Range("E12").Value = "  +1.97%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.367"
$ws.Range("E13").Value = "  -0.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "92.26"
$ws.Range("E14").Value = "  +0.64%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.519"
$ws.Range("E15").Value = "  -1.99%  "
$ws.Range("B16").Value = "TRON"
$ws.Range("C16").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.07070"
$ws.Range("E16").Value = "  -0.07%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.001"
$ws.Range("E17").Value = "  -0.13%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008726"
$ws.Range("E18").Value = "  -0.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("E19").Value = "  -0.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.73"
$ws.Range("E20").Value = "  -0.86%  "
$ws.Range("D21").Value = "26.975.29"
$ws.Range("E21").Value = "  -0.75%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.352"
$ws.Range("E22").Value = "  +0.80%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.57"
$ws.Range("E23").Value = "  -2.95%  "
$ws.Range("D24").Value = "2.092.95"
$ws.Range("E24").Value = "  +2.50%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.895"
$ws.Range("E25").Value = "  -1.91%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.89"
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.188"
$ws.Range("E27").Value = "  -1.12%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.44"
$ws.Range("E28").Value = "  -0.15%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.326"
$ws.Range("E29").Value = "  +1.30%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.64"
$ws.Range("E30").Value = "  -1.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08917"
$ws.Range("E31").Value = "  +0.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7627"
$ws.Range("E32").Value = "  -1.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.162"
$ws.Range("E33").Value = "  -0.68%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.486"
$ws.Range("E34").Value = "  -0.56%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.926"
$ws.Range("E35").Value = "  +1.18%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.000"
$ws.Range("E36").Value = "  -0.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.096"
$ws.Range("E37").Value = "  -2.73%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01960"
$ws.Range("E38").Value = "  +0.07%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05259"
$ws.Range("E39").Value = "  +0.50%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5379"
$ws.Range("E40").Value = "  +1.88%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.391"
$ws.Range("E41").Value = "  +0.39%  "
$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.927"
$ws.Range("E42").Value = "  +0.50%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.243"
$ws.Range("E43").Value = "  -0.33%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1664"
$ws.Range("E44").Value = "  -0.97%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.521"
$ws.Range("E45").Value = "  -0.82%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4959"
$ws.Range("E46").Value = "  -1.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.43"
$ws.Range("E47").Value = "  +0.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.680"
$ws.Range("E48").Value = "  +0.69%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.000"
$ws.Range("E49").Value = "  -0.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "103.31"
$ws.Range("E50").Value = "  -1.56%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06287"
$ws.Range("E51").Value = "  -0.54%  "
